$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 50 - mirrors the row layout/styling of row 49 (the last existing
# data row): column A keeps the bordered/centered "index" style, column E
# keeps the date-time number style, every other column is plain/default.
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)

$ws.Range("E49").Copy()
$ws.Range("E50").PasteSpecial(-4122)

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "denmark"
$ws.Range("C50").Value = "superliga"
$ws.Range("D50").Value = "2023-2024"
$ws.Range("E50").Value = 45191.79166666666
$ws.Range("F50").Value = "Lyngby"
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = "Vejle"
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 1.87
$ws.Range("K50").Value = "18/09/2023 08:42"
$ws.Range("L50").Value = 2.05
$ws.Range("M50").Value = "22/09/2023 18:50"
$ws.Range("N50").Value = 3.76
$ws.Range("O50").Value = "18/09/2023 08:42"
$ws.Range("P50").Value = 3.41
$ws.Range("Q50").Value = "22/09/2023 18:53"
$ws.Range("R50").Value = 3.81
$ws.Range("S50").Value = "18/09/2023 08:42"
$ws.Range("T50").Value = 3.94
$ws.Range("U50").Value = "22/09/2023 18:50"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/denmark/superliga/lyngby-vejle/6XidnLHs/"
